$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3967.375
$ws.Range("J17").Value = 3967.375
$ws.Range("L17").Value = 11902.125
$ws.Range("N17").Value = -12238.125
$ws.Range("H55").Value = 703.7778
$ws.Range("J55").Value = 934.75
$ws.Range("L55").Value = 934.75
$ws.Range("N55").Value = -1362.75
$ws.Range("H62").Value = 14038629
$ws.Range("I62").Value = 20516188
$ws.Range("J62").Value = 3916.3333
$ws.Range("K62").Value = 20516188
$ws.Range("L62").Value = 3916.3333
$ws.Range("M62").Value = -20515564
$ws.Range("N62").Value = -5164.3333
$ws.Range("H64").Value = 4655.5
$ws.Range("I64").Value = 4626
$ws.Range("J64").Value = 4803
$ws.Range("K64").Value = 4626
$ws.Range("L64").Value = 4803
$ws.Range("M64").Value = -4378
$ws.Range("N64").Value = -5299
$ws.Range("H65").Value = 14038629
$ws.Range("I65").Value = 20516188
$ws.Range("J65").Value = 3916.3333
$ws.Range("K65").Value = 102580940
$ws.Range("L65").Value = 19581.6665
$ws.Range("M65").Value = -102577820
$ws.Range("N65").Value = -25821.6665
$ws.Range("H67").Value = 4655.5
$ws.Range("I67").Value = 4626
$ws.Range("J67").Value = 4803
$ws.Range("K67").Value = 4626
$ws.Range("L67").Value = 4803
$ws.Range("M67").Value = -3768
$ws.Range("N67").Value = -6519
$ws.Range("H76").Value = 12506626
$ws.Range("I76").Value = 16674250
$ws.Range("J76").Value = 3751.5
$ws.Range("K76").Value = 16674250
$ws.Range("L76").Value = 3751.5
$ws.Range("M76").Value = -16673935
$ws.Range("N76").Value = -4381.5
$ws.Range("H79").Value = 12506626
$ws.Range("I79").Value = 16674250
$ws.Range("J79").Value = 3751.5
$ws.Range("K79").Value = 16674250
$ws.Range("L79").Value = 3751.5
$ws.Range("M79").Value = -16673158
$ws.Range("N79").Value = -5935.5
$ws.Range("H112").Value = 3498968.8
$ws.Range("J112").Value = 4134145
$ws.Range("L112").Value = 12402435
$ws.Range("N112").Value = -12404651

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 29761.637
$ws.Range("I32").Value = 29059.074
$ws.Range("J32").Value = 36787.25
$ws.Range("K32").Value = 29059.074
$ws.Range("L32").Value = 36787.25
$ws.Range("M32").Value = -28772.074
$ws.Range("N32").Value = -37361.25
$ws.Range("H61").Value = 12575.526
$ws.Range("I61").Value = 11197.1
$ws.Range("K61").Value = 11197.1
$ws.Range("M61").Value = -10985.1
$ws.Range("H63").Value = 7127.4326
$ws.Range("I63").Value = 3902.5
$ws.Range("K63").Value = 3902.5
$ws.Range("M63").Value = -3216.5
$ws.Range("H66").Value = 7127.4326
$ws.Range("I66").Value = 3902.5
$ws.Range("K66").Value = 19512.5
$ws.Range("M66").Value = -16080.5
$ws.Range("H97").Value = 1129514.1
$ws.Range("I97").Value = 1691760.2
$ws.Range("J97").Value = 5021.8184
$ws.Range("K97").Value = 1691760.2
$ws.Range("L97").Value = 5021.8184
$ws.Range("M97").Value = -1691264.2
$ws.Range("N97").Value = -6013.8184
$ws.Range("H132").Value = 7663.654
$ws.Range("I132").Value = 5779
$ws.Range("J132").Value = 10233.637
$ws.Range("K132").Value = 17337
$ws.Range("L132").Value = 30700.911
$ws.Range("M132").Value = -14807
$ws.Range("N132").Value = -35760.911
$ws.Range("H136").Value = 12575.526
$ws.Range("I136").Value = 11197.1
$ws.Range("K136").Value = 33591.3
$ws.Range("M136").Value = -31041.3
$ws.Range("H139").Value = 79602.39999999999
$ws.Range("J139").Value = 79602.39999999999
$ws.Range("L139").Value = 79602.39999999999
$ws.Range("N139").Value = -89882.39999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 40366.5
$ws.Range("J81").Value = 40366.5
$ws.Range("L81").Value = 40366.5
$ws.Range("N81").Value = -42488.5
$ws.Range("H84").Value = 40366.5
$ws.Range("J84").Value = 40366.5
$ws.Range("L84").Value = 121099.5
$ws.Range("N84").Value = -131707.5
$ws.Range("H97").Value = 30000
$ws.Range("J97").Value = 54000
$ws.Range("L97").Value = 54000
$ws.Range("N97").Value = -55982
$ws.Range("H107").Value = 1030.2273
$ws.Range("I107").Value = 792.63635
$ws.Range("K107").Value = 792.63635
$ws.Range("M107").Value = 1127.36365
$ws.Range("H134").Value = 4366.1313
$ws.Range("I134").Value = 2503.8462
$ws.Range("J134").Value = 8401.083000000001
$ws.Range("K134").Value = 7511.5386
$ws.Range("L134").Value = 25203.249
$ws.Range("M134").Value = -4976.5386
$ws.Range("N134").Value = -30273.249

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 29415836
$ws.Range("I31").Value = 52633504
$ws.Range("J31").Value = 6789.3335
$ws.Range("K31").Value = 52633504
$ws.Range("L31").Value = 6789.3335
$ws.Range("M31").Value = -52633209
$ws.Range("N31").Value = -7379.3335
$ws.Range("H34").Value = 29415836
$ws.Range("I34").Value = 52633504
$ws.Range("J34").Value = 6789.3335
$ws.Range("K34").Value = 52633504
$ws.Range("L34").Value = 6789.3335
$ws.Range("M34").Value = -52633302
$ws.Range("N34").Value = -7193.3335
$ws.Range("H51").Value = 56298
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 56298
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 56298
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = -57770
$ws.Range("H61").Value = 56298
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 56298
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 56298
$ws.Range("M61").ClearContents()
$ws.Range("N61").Value = -56994
$ws.Range("H86").Value = 6089.3335
$ws.Range("J86").Value = 7115.7144
$ws.Range("L86").Value = 7115.7144
$ws.Range("N86").Value = -9361.714400000001
$ws.Range("H89").Value = 6089.3335
$ws.Range("J89").Value = 7115.7144
$ws.Range("L89").Value = 35578.572
$ws.Range("N89").Value = -46810.572
$ws.Range("H105").Value = 2499.3333
$ws.Range("I105").Value = 1249
$ws.Range("K105").Value = 1249
$ws.Range("M105").Value = 498
$ws.Range("H107").Value = 793.2
$ws.Range("I107").Value = 835.3077
$ws.Range("K107").Value = 835.3077
$ws.Range("M107").Value = 1084.6923

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 8701.538
$ws.Range("I26").Value = 202.85715
$ws.Range("J26").Value = 18616.666
$ws.Range("K26").Value = 608.5714499999999
$ws.Range("L26").Value = 55849.99800000001
$ws.Range("M26").Value = -320.5714499999999
$ws.Range("N26").Value = -56425.99800000001
$ws.Range("H34").Value = 1089.8334
$ws.Range("I34").Value = 40.666668
$ws.Range("J34").Value = 6335.6665
$ws.Range("K34").Value = 122.000004
$ws.Range("L34").Value = 19006.9995
$ws.Range("M34").Value = -38.000004
$ws.Range("N34").Value = -19174.9995
$ws.Range("H87").Value = 3999
$ws.Range("I87").Value = 998.5
$ws.Range("J87").Value = 10000
$ws.Range("K87").Value = 2995.5
$ws.Range("L87").Value = 30000
$ws.Range("M87").Value = -1747.5
$ws.Range("N87").Value = -32496
$ws.Range("H90").Value = 3999
$ws.Range("I90").Value = 998.5
$ws.Range("J90").Value = 10000
$ws.Range("K90").Value = 8986.5
$ws.Range("L90").Value = 90000
$ws.Range("M90").Value = -2746.5
$ws.Range("N90").Value = -102480

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3252.0476
$ws.Range("I126").Value = 2315
$ws.Range("J126").Value = 4774.75
$ws.Range("K126").Value = 6945
$ws.Range("L126").Value = 14324.25
$ws.Range("M126").Value = -4475
$ws.Range("N126").Value = -19264.25
$ws.Range("H132").Value = 4462.839
$ws.Range("I132").Value = 2277.476
$ws.Range("K132").Value = 6832.428
$ws.Range("M132").Value = -4302.428
$ws.Range("H133").Value = 99999
$ws.Range("J133").Value = 99999
$ws.Range("L133").Value = 99999
$ws.Range("N133").Value = -110119
$ws.Range("H140").Value = 98333
$ws.Range("J140").Value = 98333
$ws.Range("L140").Value = 98333
$ws.Range("N140").Value = -108693

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H42").Value = 32333
$ws.Range("I42").Value = 16000
$ws.Range("K42").Value = 16000
$ws.Range("M42").Value = -15437
$ws.Range("H49").Value = 32333
$ws.Range("I49").Value = 16000
$ws.Range("K49").Value = 16000
$ws.Range("M49").Value = -15853
$ws.Range("H82").Value = 1435.8572
$ws.Range("I82").Value = 1574.875
$ws.Range("J82").Value = 1250.5
$ws.Range("K82").Value = 1574.875
$ws.Range("L82").Value = 1250.5
$ws.Range("M82").Value = -1213.875
$ws.Range("N82").Value = -1972.5
$ws.Range("H85").Value = 1435.8572
$ws.Range("I85").Value = 1574.875
$ws.Range("J85").Value = 1250.5
$ws.Range("K85").Value = 1574.875
$ws.Range("L85").Value = 1250.5
$ws.Range("M85").Value = -326.875
$ws.Range("N85").Value = -3746.5
$ws.Range("H135").Value = 95196.836
$ws.Range("J135").Value = 95196.836
$ws.Range("L135").Value = 95196.836
$ws.Range("N135").Value = -105336.836

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4333
$ws.Range("I62").Value = 3500.5
$ws.Range("J62").Value = 5998
$ws.Range("K62").Value = 3500.5
$ws.Range("L62").Value = 5998
$ws.Range("M62").Value = -2876.5
$ws.Range("N62").Value = -7246
$ws.Range("H65").Value = 4333
$ws.Range("I65").Value = 3500.5
$ws.Range("J65").Value = 5998
$ws.Range("K65").Value = 17502.5
$ws.Range("L65").Value = 29990
$ws.Range("M65").Value = -14382.5
$ws.Range("N65").Value = -36230
$ws.Range("H126").Value = 4999.0312
$ws.Range("I126").Value = 4521.7393
$ws.Range("K126").Value = 13565.2179
$ws.Range("M126").Value = -11095.2179
$ws.Range("H135").Value = 116060.57
$ws.Range("J135").Value = 116060.57
$ws.Range("L135").Value = 116060.57
$ws.Range("N135").Value = -126200.57
$ws.Range("H136").Value = 3503.6538
$ws.Range("I136").Value = 1516.5883
$ws.Range("K136").Value = 4549.7649
$ws.Range("M136").Value = -1999.7649
